# Updates the cryptocurrency price/volume table (columns D and E) on the
# active worksheet to match the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "81.521.79"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "3.162.17"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'210.36"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").Value = "'620.07"
$ws.Range("E6").Value = "  -2.55%  "
$ws.Range("D7").Value = "'0.279"
$ws.Range("E7").Value = "  +18.55%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").Value = "3.161.35"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "'0.0000250"
$ws.Range("E12").Value = "  +9.64%  "
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("E14").Value = "  -5.11%  "
$ws.Range("D15").Value = "3.741.54"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "'31.25"
$ws.Range("E16").Value = "  -1.74%  "
$ws.Range("D17").Value = "81.220.38"
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").Value = "3.160.91"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").Value = "'13.89"
$ws.Range("E20").Value = "  -4.71%  "
$ws.Range("D21").Value = "'431.32"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "'8.89"
$ws.Range("E22").Value = "  -2.89%  "
$ws.Range("D23").Value = "'5.05"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").Value = "'7.25"
$ws.Range("E24").Value = "  +5.56%  "
$ws.Range("E25").Value = "  +8.01%  "
$ws.Range("D26").Value = "3.306.26"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").Value = "'76.26"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("E28").Value = "  -4.69%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").Value = "'582.70"
$ws.Range("E31").Value = "  +10.36%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").Value = "'8.87"
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("D34").Value = "'1.49"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "'0.154"
$ws.Range("E35").Value = "  +7.36%  "
$ws.Range("E36").Value = "  +14.76%  "
$ws.Range("E37").Value = "  -1.63%  "
$ws.Range("D38").Value = "'22.65"
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "'6.07"
$ws.Range("E40").Value = "  +10.73%  "
$ws.Range("D41").Value = "'0.404"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "'2.06"
$ws.Range("E42").Value = "  +14.21%  "
$ws.Range("D43").Value = "'20.76"
$ws.Range("E43").Value = "  +3.68%  "
$ws.Range("D44").Value = "'3.00"
$ws.Range("E44").Value = "  +17.75%  "
$ws.Range("D45").Value = "'159.47"
$ws.Range("E45").Value = "  -3.21%  "
$ws.Range("D47").Value = "'186.46"
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("D48").Value = "'45.08"
$ws.Range("E48").Value = "  +4.75%  "
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("E50").Value = "  -4.76%  "
$ws.Range("D51").Value = "'25.72"
$ws.Range("E51").Value = "  -0.94%  "
